# feat: add 2022-Q1 data
#
# Inserts a new "2022-Q1" worksheet (between "2021-Q4" and "总计") holding the
# quarter's fund-holding detail rows, and updates the "总计" (totals) summary
# sheet with a new leading row for 2022-Q1 (pushing the existing 2021-Q4
# summary row down).

$wb = $excel.ActiveWorkbook
$q4 = $wb.Worksheets.Item("2021-Q4")

# Helper: write a cell as TEXT, even when the string looks like a number
# (leading-zero fund codes, "8.13"-style figures, ...) so Excel doesn't
# silently coerce it to a numeric value / drop leading zeros.
function Set-TextCell($range, [string]$text) {
    if ($text -match '^-?\d+(\.\d+)?$') {
        $range.Value = "'" + $text
    } else {
        $range.Value = $text
    }
}

# ---------------------------------------------------------------------
# 1) Insert the new "2022-Q1" sheet right after "2021-Q4"
# ---------------------------------------------------------------------
$q1 = $wb.Worksheets.Add($null, $q4)
$q1.Name = "2022-Q1"

# Match page setup / outline defaults used by the other sheets in the file
$q1.Outline.SummaryRow = 1
$q1.Outline.SummaryColumn = 1
$q1.PageSetup.LeftMargin = 54
$q1.PageSetup.RightMargin = 54
$q1.PageSetup.TopMargin = 72
$q1.PageSetup.BottomMargin = 72
$q1.PageSetup.HeaderMargin = 36
$q1.PageSetup.FooterMargin = 36

# Carry over the header-row / index-column look (bold, centered, bordered)
# from the "2021-Q4" sheet so the new sheet matches the existing style.
$q4.Range("B1:H1").Copy()
$q1.Range("B1:H1").PasteSpecial(-4122)
$q4.Range("A2:A9").Copy()
$q1.Range("A2:A9").PasteSpecial(-4122)

# Header row
Set-TextCell $q1.Range("B1") "基金代码"
Set-TextCell $q1.Range("C1") "基金名称"
Set-TextCell $q1.Range("D1") "基金规模"
Set-TextCell $q1.Range("E1") "股票总仓位"
Set-TextCell $q1.Range("F1") "仓位占比"
Set-TextCell $q1.Range("G1") "持有市值(亿元)"
Set-TextCell $q1.Range("H1") "仓位排名"

# Data rows: index, 基金代码, 基金名称, 基金规模, 股票总仓位, 仓位占比, 持有市值(亿元), 仓位排名
$rows = @(
    @(0, "010054", "万家健康产业混合A",                     "8.13", "86.63", "4.01", "0.3260", 3),
    @(1, "010055", "万家健康产业混合C",                     "3.36", "86.63", "4.01", "0.1347", 3),
    @(2, "005108", "圆信永丰双利优选定期开放灵活配置混合",   "1.89", "94.60", "5.43", "0.1026", 6),
    @(3, "009893", "摩根士丹利华鑫优悦安和混合",             "0.91", "93.90", "7.09", "0.0645", 8),
    @(4, "160921", "大成多策略混合(LOF)",                    "1.13", "79.19", "3.91", "0.0442", 8),
    @(5, "001965", "圆信永丰兴源灵活配置混合A",              "0.76", "93.43", "4.52", "0.0344", 9),
    @(6, "001966", "圆信永丰兴源灵活配置混合C",              "0.25", "93.43", "4.52", "0.0113", 9),
    @(7, "006274", "圆信永丰医药健康混合",                   "0.18", "93.66", "5.26", "0.0095", 6)
)

$r = 2
foreach ($row in $rows) {
    $q1.Range("A$r").Value = $row[0]
    Set-TextCell $q1.Range("B$r") $row[1]
    Set-TextCell $q1.Range("C$r") $row[2]
    Set-TextCell $q1.Range("D$r") $row[3]
    Set-TextCell $q1.Range("E$r") $row[4]
    Set-TextCell $q1.Range("F$r") $row[5]
    Set-TextCell $q1.Range("G$r") $row[6]
    $q1.Range("H$r").Value = $row[7]
    $r = $r + 1
}

# ---------------------------------------------------------------------
# 2) Update the "总计" sheet: push the existing 2021-Q4 row down and add
#    a new 2022-Q1 row above it (row 2).
#
# NB: fetch this handle *after* inserting the new sheet above - this host's
# Worksheets.Item(name) binds positionally, so a handle grabbed before the
# insert would silently start tracking the new sheet once indices shift.
# ---------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")
$total.Rows.Item(2).Insert()

# The insert carries the row-1 formatting down into the new row 2; strip it
# back to the unstyled look the data rows use, then restore the index-
# column (A) style by copying it from row 3 (the shifted-down 2021-Q4 row).
$total.Range("A2:D2").ClearFormats()
$total.Range("A3").Copy()
$total.Range("A2").PasteSpecial(-4122)

$total.Range("A2").Value = 0
Set-TextCell $total.Range("B2") "2022-Q1"
$total.Range("C2").Value = 8
$total.Range("D2").Value = 0.73

# Restore the original 2021-Q4 row's values explicitly (row insert can
# otherwise leave the shifted row's cells blank/stale in some hosts).
$total.Range("A3").Value = 1
Set-TextCell $total.Range("B3") "2021-Q4"
$total.Range("C3").Value = 13
$total.Range("D3").Value = 5.47

$q4.Activate()
